$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "UEprobInd" regression (5th column, "spending IIIII") is dropped from
# the table entirely: delete its whole column (F). This shifts the last
# regression column ("UEprobAgg" / old "spending IIIIII", column G) left
# into the F position, carrying its coefficient/se/N values with it.
$ws.Range("F1:F15").EntireColumn.Delete()

# The "UEprobInd" row pair (row label + blank se row) is also dropped from
# the left-hand label column: delete rows 10:11, shifting UEprobAgg's row
# and the trailing N / R2 rows up by two.
$ws.Range("A10:A11").EntireRow.Delete()

# The column that slid into F used to be headed "spending IIIIII"; relabel
# it to match its new position as the 5th (and now last) regression,
# "spending IIIII".
$ws.Range("F1").Value = "spending IIIII"
